# Update search term values in column A (rows 1-100) per KW-RANK diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    'extra long basketball shorts for men',
    'compression pants running',
    'boys baseball compression sleeve',
    'stretch mark for men',
    'compression shorts men pack',
    'knee pads youth girls',
    'knee sleeves basketball',
    'xl youth baseball pants',
    'biking knee sleeve',
    'boy baseball pants',
    'men black compression pants',
    'knee pads exercise',
    'youth basketball compression shorts',
    'hockey pads',
    'calf strain compression sleeve',
    'mens degree sport',
    'men long shorts below the knee',
    'knee pads construction gel',
    'work knee pads construction',
    'padded football sleeve',
    'spandex tights',
    'adult black baseball pants',
    'girls baseball pants',
    'mens small black baseball pants',
    'youth softball compression sleeve',
    'padded shorts football',
    '3/4 shorts for men',
    'silicone strips for stretch marks',
    'calf compression sleeve basketball',
    'gel knee pads for construction',
    'youth baseball short pants',
    'exercise knee pad',
    'girls basketball shorts size 6',
    'mens short tights',
    'boys knee length shorts',
    'rodilleras de volleyball',
    'rodilleras volleyball',
    'compression running pants men',
    'knee pads cheap',
    'exercise knee pads',
    'pad knee',
    'compression workout pants men',
    'compression sleeve for bursitis',
    'mens running tights shorts',
    'mens compression girdle',
    'knee pads outdoor',
    'motorcycle knee pads for men',
    'water pants men',
    'pant baseball men',
    'volleyball knee pad',
    'knee pads for volleyball',
    'pro tights men',
    'sleeve knee pads',
    'womens compression leggings',
    'airsoft knee pads',
    'skateboard knee pads',
    'starter youth compression pants',
    'trolls knee pads',
    'mens compression pants marvel',
    'mouthguard basketball youth',
    'mcdavid youth knee pads',
    'nike compression leggings',
    'jordan flight mens basketball pants',
    'athletic compression pants',
    'nike pro compression leggings men',
    'emoji knee pads',
    'mens compression tights 3 4',
    'mens compression tights nike',
    'mens compression tights white',
    'eastbay compression pants',
    'women compression leggings',
    'ladies compression pants',
    'womans compression leggings',
    'mcdavid compression pants',
    'knee brace basketball youth',
    'basketball kids knee pads',
    'basketball knee pads mcdavid',
    'basketball youth jersey',
    'morris compression knee pads',
    'wonens compression leggings',
    'youth baketball knee pads',
    'olympic mens basketball',
    'track leggings men',
    'nike pro dry mens basketball tights',
    'goalie knee protectors',
    'basketball clothes for men',
    'mens running thermal compression pants',
    'capri pants for men adidas',
    'men running tights',
    'mens running tight',
    'men tights nike',
    'men tights short',
    'mens tights dance',
    'men running tight',
    'men tights green',
    'men tights pack',
    'mens tights grey',
    'mens tights nike',
    'mens tights pink',
    'men legging nike'
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $newValues[$i]
}
